$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 23: the email column was pointing at the wrong (duplicate)
# shared string; it should hold Vishal M's actual email address. This
# also introduces a new shared-string entry, matching the sharedStrings.xml
# insertion in the target diff.
$ws.Range("A23").Value = "vishal.maheshab@gmail.com"

# --- Column A is now wide enough to fit full e-mail addresses.
$ws.Range("A1").EntireColumn.ColumnWidth = 40.57642857142857

# --- Data rows 2-28 grow slightly taller (18.75pt -> 19.5pt).
$ws.Rows("2:28").RowHeight = 19.5

# --- The header row (A1:M1) and the numeric data block (D2:M35) share
# the same underlying font, which switches from a theme-based color to
# an explicit black RGB color (the duplicate font used only by D2:M35
# is dropped in favor of that shared font).
$ws.Range("A1:M1").Font.Color = 0
$ws.Range("D2:M35").Font.Color = 0
